# Update "想去人数" (column F) figures across three worksheets to match
# the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 37161
$ws.Range("F5").Value  = 763
$ws.Range("F6").Value  = 473
$ws.Range("F8").Value  = 463
$ws.Range("F9").Value  = 833
$ws.Range("F11").Value = 685
$ws.Range("F14").Value = 629
$ws.Range("F16").Value = 463
$ws.Range("F17").Value = 439
$ws.Range("F18").Value = 1152
$ws.Range("F20").Value = 801
$ws.Range("F21").Value = 2481
$ws.Range("F22").Value = 982
$ws.Range("F24").Value = 100
$ws.Range("F27").Value = 751
$ws.Range("F28").Value = 48

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 328

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 37161
$ws.Range("F6").Value  = 763
$ws.Range("F7").Value  = 473
$ws.Range("F10").Value = 463
$ws.Range("F12").Value = 328
$ws.Range("F14").Value = 833
$ws.Range("F16").Value = 685
$ws.Range("F24").Value = 629
$ws.Range("F26").Value = 463
$ws.Range("F27").Value = 439
$ws.Range("F28").Value = 1152
$ws.Range("F30").Value = 801
$ws.Range("F31").Value = 2481
$ws.Range("F32").Value = 982
$ws.Range("F34").Value = 100
$ws.Range("F38").Value = 751
$ws.Range("F39").Value = 48

$wb.Save()
